$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ax..gz (C:H) values for existing rows 2-21 with newly recorded sensor data
$ws.Range("C2").Value2 = -1.610764980316162
$ws.Range("D2").Value2 = 1.577288150787354
$ws.Range("E2").Value2 = 0.315173327922821
$ws.Range("F2").Value2 = 0.007011067026972018
$ws.Range("G2").Value2 = -0.03245915641838846
$ws.Range("H2").Value2 = 0.03166781107641074
$ws.Range("C3").Value2 = -1.949196338653564
$ws.Range("D3").Value2 = 1.61798095703125
$ws.Range("E3").Value2 = 0.4421060681343078
$ws.Range("F3").Value2 = -0.0004164990228177805
$ws.Range("G3").Value2 = -0.01167585594918232
$ws.Range("H3").Value2 = -0.0006108652715655517
$ws.Range("C4").Value2 = -1.840867042541504
$ws.Range("D4").Value2 = 1.58759355545044
$ws.Range("E4").Value2 = 0.5192338824272156
$ws.Range("F4").Value2 = -0.0122173046693205
$ws.Range("G4").Value2 = 0.0074830991216003
$ws.Range("H4").Value2 = 0.0039706239476799
$ws.Range("C5").Value2 = -1.788941383361816
$ws.Range("D5").Value2 = 1.582527160644531
$ws.Range("E5").Value2 = 0.4804926216602325
$ws.Range("F5").Value2 = -0.04220523541285226
$ws.Range("G5").Value2 = 0.007899598006836298
$ws.Range("H5").Value2 = 0.001610462892461813
$ws.Range("C6").Value2 = -1.827802658081055
$ws.Range("D6").Value2 = 1.591060638427734
$ws.Range("E6").Value2 = 0.4033206701278686
$ws.Range("F6").Value2 = -0.01571589649062261
$ws.Range("G6").Value2 = -0.009926560250195525
$ws.Range("H6").Value2 = 0.002873843345283052
$ws.Range("C7").Value2 = -1.872776985168457
$ws.Range("D7").Value2 = 1.712420463562011
$ws.Range("E7").Value2 = 0.4200127720832824
$ws.Range("F7").Value2 = 0.007719115523452058
$ws.Range("G7").Value2 = -0.02122756669467137
$ws.Range("H7").Value2 = 0.02797485068863757
$ws.Range("C8").Value2 = -1.756282329559326
$ws.Range("D8").Value2 = 1.523788452148438
$ws.Range("E8").Value2 = 0.3277221620082855
$ws.Range("F8").Value2 = -0.01150925694541499
$ws.Range("G8").Value2 = -0.008524346549000439
$ws.Range("H8").Value2 = 0.01102334066209471
$ws.Range("C9").Value2 = -1.891244411468506
$ws.Range("D9").Value2 = 1.559478759765625
$ws.Range("E9").Value2 = 0.1866782307624817
$ws.Range("F9").Value2 = -0.01870080676268442
$ws.Range("G9").Value2 = -0.01660442801816813
$ws.Range("H9").Value2 = -0.0122173046693205
$ws.Range("C10").Value2 = -1.622483730316162
$ws.Range("D10").Value2 = 1.571603775024414
$ws.Range("E10").Value2 = 0.2709611356258392
$ws.Range("F10").Value2 = -0.01634064570746632
$ws.Range("G10").Value2 = -0.02958531457592136
$ws.Range("H10").Value2 = -0.03082092817534093
$ws.Range("C11").Value2 = -1.867420673370361
$ws.Range("D11").Value2 = 1.565328121185303
$ws.Range("E11").Value2 = 0.3432579040527344
$ws.Range("F11").Value2 = -0.01731247691945591
$ws.Range("G11").Value2 = -0.0009301814504645012
$ws.Range("H11").Value2 = -0.02040845257314771
$ws.Range("C12").Value2 = -1.858330726623535
$ws.Range("D12").Value2 = 1.508580207824707
$ws.Range("E12").Value2 = 0.2882210314273834
$ws.Range("F12").Value2 = -0.03719336404041793
$ws.Range("G12").Value2 = 0.1217704361135306
$ws.Range("H12").Value2 = 0.03431951999664297
$ws.Range("C13").Value2 = -1.921436786651612
$ws.Range("D13").Value2 = 1.638089656829834
$ws.Range("E13").Value2 = 0.3902863562107086
$ws.Range("F13").Value2 = -0.1360424668951469
$ws.Range("G13").Value2 = 0.3186911859295585
$ws.Range("H13").Value2 = 0.05561650341207323
$ws.Range("C14").Value2 = -1.84520435333252
$ws.Range("D14").Value2 = 1.668948650360107
$ws.Range("E14").Value2 = 0.1994338035583496
$ws.Range("F14").Value2 = -0.2708076590841467
$ws.Range("G14").Value2 = 0.3116106986999512
$ws.Range("H14").Value2 = 0.0008329986171288634
$ws.Range("C15").Value2 = -1.829729080200196
$ws.Range("D15").Value2 = 1.47260046005249
$ws.Range("E15").Value2 = 0.5425764322280884
$ws.Range("F15").Value2 = -0.3139847218990326
$ws.Range("G15").Value2 = 0.365755558013916
$ws.Range("H15").Value2 = -0.0809396430850029
$ws.Range("C16").Value2 = -1.84237813949585
$ws.Range("D16").Value2 = 1.411740303039551
$ws.Range("E16").Value2 = 1.237700462341309
$ws.Range("F16").Value2 = -0.09199075129899148
$ws.Range("G16").Value2 = 0.302864204753529
$ws.Range("H16").Value2 = -0.2070000191981139
$ws.Range("C17").Value2 = -2.159674644470215
$ws.Range("D17").Value2 = 1.136321544647217
$ws.Range("E17").Value2 = 3.083052396774292
$ws.Range("F17").Value2 = 0.01230060461569882
$ws.Range("G17").Value2 = -0.3993948210369431
$ws.Range("H17").Value2 = -0.9949190426956491
$ws.Range("C18").Value2 = -1.937846660614014
$ws.Range("D18").Value2 = 1.245656490325928
$ws.Range("E18").Value2 = 3.884034872055054
$ws.Range("F18").Value2 = 0.3243139045820992
$ws.Range("G18").Value2 = -0.5972596244378523
$ws.Range("H18").Value2 = -1.587486072020097
$ws.Range("C19").Value2 = -3.061091899871826
$ws.Range("D19").Value2 = 0.4542388916015625
$ws.Range("E19").Value2 = 2.74072003364563
$ws.Range("F19").Value2 = -0.9380390860817602
$ws.Range("G19").Value2 = -1.513376929543232
$ws.Range("H19").Value2 = -3.357259793715038
$ws.Range("C20").Value2 = -2.979402542114258
$ws.Range("D20").Value2 = -0.1682605743408203
$ws.Range("E20").Value2 = 2.131577730178833
$ws.Range("F20").Value2 = -0.2036124901338088
$ws.Range("G20").Value2 = -0.2562440525401761
$ws.Range("H20").Value2 = -0.1245054331692872
$ws.Range("C21").Value2 = -3.620201349258423
$ws.Range("D21").Value2 = -1.011659622192383
$ws.Range("E21").Value2 = 2.961694955825806
$ws.Range("F21").Value2 = 0.2954227382486529
$ws.Range("G21").Value2 = 0.008163343776357967
$ws.Range("H21").Value2 = 0.4832359877499712

# Append 10 new rows (22-31) of data
$ws.Range("A22").Value2 = 2000
$ws.Range("B22").Value2 = "falling"
$ws.Range("C22").Value2 = -5.750590324401856
$ws.Range("D22").Value2 = -1.630428791046143
$ws.Range("E22").Value2 = 4.153227806091309
$ws.Range("F22").Value2 = -0.3702815256335535
$ws.Range("G22").Value2 = -0.7731332995674872
$ws.Range("H22").Value2 = -1.987241875041621
$ws.Range("A23").Value2 = 2100
$ws.Range("B23").Value2 = "falling"
$ws.Range("C23").Value2 = -4.631756782531738
$ws.Range("D23").Value2 = 0.1280508041381836
$ws.Range("E23").Value2 = 0.0708565711975097
$ws.Range("F23").Value2 = -0.1218537308953028
$ws.Range("G23").Value2 = -0.1149676279588187
$ws.Range("H23").Value2 = -0.5308419032530356
$ws.Range("A24").Value2 = 2200
$ws.Range("B24").Value2 = "falling"
$ws.Range("C24").Value2 = -6.701959609985352
$ws.Range("D24").Value2 = 13.03308868408203
$ws.Range("E24").Value2 = -3.955403804779053
$ws.Range("F24").Value2 = -0.7092561884359876
$ws.Range("G24").Value2 = -1.314401550726458
$ws.Range("H24").Value2 = -0.9030254652554355
$ws.Range("A25").Value2 = 2300
$ws.Range("B25").Value2 = "falling"
$ws.Range("C25").Value2 = -9.679704666137695
$ws.Range("D25").Value2 = 1.827943325042725
$ws.Range("E25").Value2 = 0.2789157629013061
$ws.Range("F25").Value2 = -0.9390387182885965
$ws.Range("G25").Value2 = -0.6636356765573627
$ws.Range("H25").Value2 = 0.2282692268490727
$ws.Range("A26").Value2 = 2400
$ws.Range("B26").Value2 = "falling"
$ws.Range("C26").Value2 = 3.086193084716797
$ws.Range("D26").Value2 = -2.266827583312988
$ws.Range("E26").Value2 = 2.888416528701782
$ws.Range("F26").Value2 = -0.3692680299282074
$ws.Range("G26").Value2 = -2.486374378204346
$ws.Range("H26").Value2 = 2.66978645324707
$ws.Range("A27").Value2 = 2500
$ws.Range("B27").Value2 = "falling"
$ws.Range("C27").Value2 = 1.258205413818359
$ws.Range("D27").Value2 = 1.049188137054443
$ws.Range("E27").Value2 = -6.970268726348877
$ws.Range("F27").Value2 = 0.7195992388508516
$ws.Range("G27").Value2 = 0.8889339186928424
$ws.Range("H27").Value2 = 0.3780700483105441
$ws.Range("A28").Value2 = 2600
$ws.Range("B28").Value2 = "falling"
$ws.Range("C28").Value2 = -2.633459091186523
$ws.Range("D28").Value2 = 5.726268291473389
$ws.Range("E28").Value2 = 6.590025901794434
$ws.Range("F28").Value2 = -1.744853258132938
$ws.Range("G28").Value2 = -2.628608725287703
$ws.Range("H28").Value2 = 2.057380207560282
$ws.Range("A29").Value2 = 2700
$ws.Range("B29").Value2 = "falling"
$ws.Range("C29").Value2 = 4.84531307220459
$ws.Range("D29").Value2 = -0.8891797065734863
$ws.Range("E29").Value2 = -5.761303424835205
$ws.Range("F29").Value2 = -1.074428666721691
$ws.Range("G29").Value2 = -0.7690515951676815
$ws.Range("H29").Value2 = 1.639270782470704
$ws.Range("A30").Value2 = 2800
$ws.Range("B30").Value2 = "falling"
$ws.Range("C30").Value2 = 0.9592771530151368
$ws.Range("D30").Value2 = 3.222768306732178
$ws.Range("E30").Value2 = 1.57231593132019
$ws.Range("F30").Value2 = -0.2539533461359424
$ws.Range("G30").Value2 = -0.02276861396702762
$ws.Range("H30").Value2 = 0.6248873986981154
$ws.Range("A31").Value2 = 2900
$ws.Range("B31").Value2 = "falling"
$ws.Range("C31").Value2 = 1.009872436523438
$ws.Range("D31").Value2 = -0.1888983249664306
$ws.Range("E31").Value2 = 0.460105299949646
$ws.Range("F31").Value2 = -0.05406157367608758
$ws.Range("G31").Value2 = 0.4298825345256099
$ws.Range("H31").Value2 = -0.1093171049248081
